$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.915.78'
$ws.Range("E2").Value = '  -0.67%  '

$ws.Range("D3").Value = '1.640.79'
$ws.Range("E3").Value = '  -0.32%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5039'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2568'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06383'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.57'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07769'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.29%  '

$ws.Range("D12").Value = '1.659.06'
$ws.Range("E12").Value = '  +0.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.267'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.21%  '

$ws.Range("D14").Value = '1.868.09'
$ws.Range("E14").Value = '  -0.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5416'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.90%  '

$ws.Range("D16").Value = '0.0₅7855'
$ws.Range("E16").Value = '  -1.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.64'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.32%  '

$ws.Range("D18").Value = '25.954.82'
$ws.Range("E18").Value = '  -0.58%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '198.05'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.377'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.938'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.88%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.971'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.006'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.82%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.872'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.40%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.82'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.69%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1142'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.838'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.240'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04897'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.251'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.186'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.529'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.364'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8911'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.95%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.601'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.98%  '

$ws.Range("D38").Value = '1.135.35'
$ws.Range("E38").Value = '  -1.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5541'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.96%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01563'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.002'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.12%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.683'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8162'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.59%  '

$ws.Range("D45").Value = '0.0₈119'
$ws.Range("E45").Value = '  +4.10%  '

$ws.Range("D46").Value = '1.778.19'
$ws.Range("E46").Value = '  -0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4511'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.008'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.18'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.18%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05081'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.007'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.63%  '

